$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text to preserve exact formatting
# (leading/trailing zeros, multi-dot values, etc.), matching the original inlineStr cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.157.67"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.82"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.78"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9940"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3931"
$ws.Range("E7").Value = "  +3.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3506"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.34"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.204"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9929"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.562"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.821.90"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.231"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06706"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.63"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9940"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.00"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.605"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.179.66"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.85"
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.408"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.543"
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.55"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "155.29"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.026.13"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.28"
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.224"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.031"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08869"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.41"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.570"
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6977"
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02443"
$ws.Range("E38").Value = "  +4.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06584"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("E40").Value = "  -4.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2233"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.272"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.599"
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.74"
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6511"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.881"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.179"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.49"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07234"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "80.85"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.160"
$ws.Range("E51").Value = "  +3.29%  "
